$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(53, 1).Copy($ws.Cells.Item(54, 1))
$ws.Cells.Item(54, 1).Value = 45986

$ws.Cells.Item(54, 2).Value = 2025
$ws.Cells.Item(54, 3).Value = 1.049317648994741
$ws.Cells.Item(54, 4).Value = 2026
$ws.Cells.Item(54, 5).Value = 0.07146359800258573
